$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "time_taken" in F1, copying header style from E1
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("F1").Value = "time_taken"

# Fill time_taken values for each data row (2-45)
$ws.Range("F2").Value = "2021-10-05 13:38:36.552832"
$ws.Range("F3").Value = "2021-10-05 13:38:36.552845"
$ws.Range("F4").Value = "2021-10-05 13:38:36.552849"
$ws.Range("F5").Value = "2021-10-05 13:38:36.552852"
$ws.Range("F6").Value = "2021-10-05 13:38:36.552856"
$ws.Range("F7").Value = "2021-10-05 13:38:36.552859"
$ws.Range("F8").Value = "2021-10-05 13:38:36.552862"
$ws.Range("F9").Value = "2021-10-05 13:38:36.552865"
$ws.Range("F10").Value = "2021-10-05 13:38:36.552868"
$ws.Range("F11").Value = "2021-10-05 13:38:36.552871"
$ws.Range("F12").Value = "2021-10-05 13:38:36.552874"
$ws.Range("F13").Value = "2021-10-05 13:38:36.552877"
$ws.Range("F14").Value = "2021-10-05 13:38:36.552880"
$ws.Range("F15").Value = "2021-10-05 13:38:36.552883"
$ws.Range("F16").Value = "2021-10-05 13:38:36.552886"
$ws.Range("F17").Value = "2021-10-05 13:38:36.552889"
$ws.Range("F18").Value = "2021-10-05 13:38:36.552892"
$ws.Range("F19").Value = "2021-10-05 13:38:36.552895"
$ws.Range("F20").Value = "2021-10-05 13:38:36.552898"
$ws.Range("F21").Value = "2021-10-05 13:38:36.552901"
$ws.Range("F22").Value = "2021-10-05 13:38:36.552904"
$ws.Range("F23").Value = "2021-10-05 13:38:36.552908"
$ws.Range("F24").Value = "2021-10-05 13:38:36.552911"
$ws.Range("F25").Value = "2021-10-05 13:38:36.552914"
$ws.Range("F26").Value = "2021-10-05 13:38:36.552917"
$ws.Range("F27").Value = "2021-10-05 13:38:36.552921"
$ws.Range("F28").Value = "2021-10-05 13:38:36.552924"
$ws.Range("F29").Value = "2021-10-05 13:38:36.552927"
$ws.Range("F30").Value = "2021-10-05 13:38:36.552930"
$ws.Range("F31").Value = "2021-10-05 13:38:36.552933"
$ws.Range("F32").Value = "2021-10-05 13:38:36.552936"
$ws.Range("F33").Value = "2021-10-05 13:38:36.552938"
$ws.Range("F34").Value = "2021-10-05 13:38:36.552942"
$ws.Range("F35").Value = "2021-10-05 13:38:36.552945"
$ws.Range("F36").Value = "2021-10-05 13:38:36.552948"
$ws.Range("F37").Value = "2021-10-05 13:38:36.552951"
$ws.Range("F38").Value = "2021-10-05 13:38:36.552954"
$ws.Range("F39").Value = "2021-10-05 13:38:36.552957"
$ws.Range("F40").Value = "2021-10-05 13:38:36.552960"
$ws.Range("F41").Value = "2021-10-05 13:38:36.552963"
$ws.Range("F42").Value = "2021-10-05 13:38:36.552966"
$ws.Range("F43").Value = "2021-10-05 13:38:36.552969"
$ws.Range("F44").Value = "2021-10-05 13:38:36.552972"
$ws.Range("F45").Value = "2021-10-05 13:38:36.552975"
